$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old I1 header cell ("Ethnicity") so the header row ends at H1
$ws.Range("I1").ClearContents()

# Fill in row 2 with the new case data
$ws.Range("A2").Value = "CTDC-46730"
$ws.Range("B2").Value = "NCI-MATCH"
$ws.Range("C2").Value = "Q"
$ws.Range("D2").Value = "Ado-trastuzumab Emtansine"
$ws.Range("E2").Value = "Adenocarcinoma of the cervix"
$ws.Range("F2").Value = "FEMALE"
$ws.Range("G2").Value = "UNKNOWN"
$ws.Range("H2").Value = "UNKNOWN"
